$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy styles from existing rows so new cells match the formatting (B: date style, D: centered style)
$ws.Range("B577").Copy() | Out-Null
$ws.Range("B578:B585").PasteSpecial(-4122) | Out-Null
$ws.Range("D566").Copy() | Out-Null
$ws.Range("D578:D585").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 578
$ws.Range("A578").Value = "Entrainement"
$ws.Range("B578").Value = 45930
$ws.Range("C578").Value = "Global"
$ws.Range("D578").Value = "J+3"
$ws.Range("E578").Value = "Kamal Bafounta"
$ws.Range("F578").Value = "center midfield"
$ws.Range("G578").Value = "01:33:33"
$ws.Range("H578").Value = 7.06
$ws.Range("I578").Value = 0.14
$ws.Range("J578").Value = 6.91
$ws.Range("K578").Value = 0.15
$ws.Range("L578").Value = 0
$ws.Range("M578").Value = 0
$ws.Range("N578").Value = 0
$ws.Range("O578").Value = 0
$ws.Range("P578").Value = 4.45
$ws.Range("Q578").Value = 20
$ws.Range("R578").Value = 4.4
$ws.Range("S578").Value = 31
$ws.Range("T578").Value = 1
$ws.Range("U578").Value = 24
$ws.Range("V578").Value = 5

# Row 579
$ws.Range("A579").Value = "Entrainement"
$ws.Range("B579").Value = 45930
$ws.Range("C579").Value = "Global"
$ws.Range("D579").Value = "J+3"
$ws.Range("E579").Value = "Karim Belmahi"
$ws.Range("F579").Value = "left forward"
$ws.Range("G579").Value = "01:32:46"
$ws.Range("H579").Value = 7.31
$ws.Range("I579").Value = 0.24
$ws.Range("J579").Value = 7.07
$ws.Range("K579").Value = 0.23
$ws.Range("L579").Value = 0.01
$ws.Range("M579").Value = 0
$ws.Range("N579").Value = 0
$ws.Range("O579").Value = 0
$ws.Range("P579").Value = 4.64
$ws.Range("Q579").Value = 22.21
$ws.Range("R579").Value = 4.23
$ws.Range("S579").Value = 37
$ws.Range("T579").Value = 2
$ws.Range("U579").Value = 29
$ws.Range("V579").Value = 1

# Row 580
$ws.Range("A580").Value = "Entrainement"
$ws.Range("B580").Value = 45930
$ws.Range("C580").Value = "Global"
$ws.Range("D580").Value = "J+3"
$ws.Range("E580").Value = "Karahali Souaré"
$ws.Range("F580").Value = "right forward"
$ws.Range("G580").Value = "01:32:46"
$ws.Range("H580").Value = 7.37
$ws.Range("I580").Value = 0.2
$ws.Range("J580").Value = 7.15
$ws.Range("K580").Value = 0.19
$ws.Range("L580").Value = 0.02
$ws.Range("M580").Value = 0
$ws.Range("N580").Value = 0
$ws.Range("O580").Value = 1
$ws.Range("P580").Value = 4.35
$ws.Range("Q580").Value = 25.53
$ws.Range("R580").Value = 5.67
$ws.Range("S580").Value = 80
$ws.Range("T580").Value = 20
$ws.Range("U580").Value = 57
$ws.Range("V580").Value = 23

# Row 581
$ws.Range("A581").Value = "Entrainement"
$ws.Range("B581").Value = 45930
$ws.Range("C581").Value = "Global"
$ws.Range("D581").Value = "J+3"
$ws.Range("E581").Value = "Omar Benyounes"
$ws.Range("F581").Value = "center midfield"
$ws.Range("G581").Value = "01:33:05"
$ws.Range("H581").Value = 7.21
$ws.Range("I581").Value = 0.25
$ws.Range("J581").Value = 6.94
$ws.Range("K581").Value = 0.23
$ws.Range("L581").Value = 0.03
$ws.Range("M581").Value = 0
$ws.Range("N581").Value = 0
$ws.Range("O581").Value = 1
$ws.Range("P581").Value = 4.53
$ws.Range("Q581").Value = 25.62
$ws.Range("R581").Value = 4.63
$ws.Range("S581").Value = 43
$ws.Range("T581").Value = 2
$ws.Range("U581").Value = 25
$ws.Range("V581").Value = 5

# Row 582
$ws.Range("A582").Value = "Entrainement"
$ws.Range("B582").Value = 45930
$ws.Range("C582").Value = "Global"
$ws.Range("D582").Value = "J+3"
$ws.Range("E582").Value = "Malik Boussaid"
$ws.Range("F582").Value = "right back"
$ws.Range("G582").Value = "01:33:42"
$ws.Range("H582").Value = 7.92
$ws.Range("I582").Value = 0.24
$ws.Range("J582").Value = 7.67
$ws.Range("K582").Value = 0.22
$ws.Range("L582").Value = 0.03
$ws.Range("M582").Value = 0
$ws.Range("N582").Value = 0
$ws.Range("O582").Value = 0
$ws.Range("P582").Value = 4.47
$ws.Range("Q582").Value = 23.99
$ws.Range("R582").Value = 4.37
$ws.Range("S582").Value = 39
$ws.Range("T582").Value = 8
$ws.Range("U582").Value = 41
$ws.Range("V582").Value = 11

# Row 583
$ws.Range("A583").Value = "Entrainement"
$ws.Range("B583").Value = 45930
$ws.Range("C583").Value = "Global"
$ws.Range("D583").Value = "J+3"
$ws.Range("E583").Value = "Mattheo Haon"
$ws.Range("F583").Value = "right back"
$ws.Range("G583").Value = "01:29:59"
$ws.Range("H583").Value = 7.31
$ws.Range("I583").Value = 0.32
$ws.Range("J583").Value = 6.99
$ws.Range("K583").Value = 0.26
$ws.Range("L583").Value = 0.06
$ws.Range("M583").Value = 0
$ws.Range("N583").Value = 0
$ws.Range("O583").Value = 0
$ws.Range("P583").Value = 4.78
$ws.Range("Q583").Value = 24.38
$ws.Range("R583").Value = 4.52
$ws.Range("S583").Value = 44
$ws.Range("T583").Value = 9
$ws.Range("U583").Value = 39
$ws.Range("V583").Value = 7

# Row 584
$ws.Range("A584").Value = "Entrainement"
$ws.Range("B584").Value = 45930
$ws.Range("C584").Value = "Global"
$ws.Range("D584").Value = "J+3"
$ws.Range("E584").Value = "Amine Taiar"
$ws.Range("F584").Value = "center back"
$ws.Range("G584").Value = "01:32:46"
$ws.Range("H584").Value = 7.39
$ws.Range("I584").Value = 0.22
$ws.Range("J584").Value = 7.16
$ws.Range("K584").Value = 0.17
$ws.Range("L584").Value = 0.06
$ws.Range("M584").Value = 0
$ws.Range("N584").Value = 0
$ws.Range("O584").Value = 0
$ws.Range("P584").Value = 4.05
$ws.Range("Q584").Value = 22.82
$ws.Range("R584").Value = 6.78
$ws.Range("S584").Value = 39
$ws.Range("T584").Value = 8
$ws.Range("U584").Value = 42
$ws.Range("V584").Value = 13

# Row 585
$ws.Range("A585").Value = "Entrainement"
$ws.Range("B585").Value = 45930
$ws.Range("C585").Value = "Global"
$ws.Range("D585").Value = "J+3"
$ws.Range("E585").Value = "Emmanuel Valey"
$ws.Range("F585").Value = "left forward"
$ws.Range("G585").Value = "01:31:50"
$ws.Range("H585").Value = 7.51
$ws.Range("I585").Value = 0.12
$ws.Range("J585").Value = 7.38
$ws.Range("K585").Value = 0.12
$ws.Range("L585").Value = 0.01
$ws.Range("M585").Value = 0
$ws.Range("N585").Value = 0
$ws.Range("O585").Value = 0
$ws.Range("P585").Value = 4.13
$ws.Range("Q585").Value = 23.19
$ws.Range("R585").Value = 4.41
$ws.Range("S585").Value = 48
$ws.Range("T585").Value = 6
$ws.Range("U585").Value = 59
$ws.Range("V585").Value = 13

$ws.Range("A578:A585").Select() | Out-Null